# "Zweiter CI mit div Aenderungen" - split the lone paragraph in two:
# keep the original sentence in its own paragraph, add a new paragraph
# with the follow-up sentence, a bold/enlarged "Fett " (still wrapped by
# the pre-existing _GoBack bookmark) and a closing "markiert!!", then
# append one more empty trailing paragraph.

$d = $word.ActiveDocument

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Grab paragraph 1's own opening <w:p ...> tag (rsid/paraId attributes)
# so we can reproduce it verbatim instead of guessing its attributes.
$p1 = $d.Paragraphs(1).Range
$p1Oxml = $p1.WordOpenXML
if ($p1Oxml -match '(<w:p [^>]*>)') {
    $p1Open = $matches[1]
} else {
    $p1Open = "<w:p $wNs>"
}

$firstSentence = "Das ist der erste Check in als TEST!! "
$secondSentence = "Nun habe ich einen Text eingefügt und danach noch ein Wort "

$newXml = $p1Open + "<w:r><w:t xml:space='preserve'>$firstSentence</w:t></w:r></w:p>" +
    "<w:p $wNs>" +
        "<w:r><w:t xml:space='preserve'>$secondSentence</w:t></w:r>" +
        "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
        "<w:r><w:rPr><w:b/><w:sz w:val='72'/></w:rPr><w:t>Fett</w:t></w:r>" +
        "<w:r><w:rPr><w:sz w:val='72'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
        "<w:bookmarkEnd w:id='0'/>" +
        "<w:r><w:t>markiert!!</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs/>"

# Replace the whole (only) paragraph with the two reconstructed
# paragraphs plus a trailing empty one.
[void]$p1.InsertXML($newXml)
